# ---------------------------------------------------------------------------
# "added road map slide"
#
# 1) Slide 1 ("Southern Spectroscopic Survey Instrument"): split the run
#    "capabilities matched to LSST and CMB-S4 survey areas and depths" so
#    that "matched to LSST and CMB-S4 survey areas and depths" becomes
#    underlined.
# 2) Add a new slide 2 ("Road Map for Spectroscopy") right after slide 1,
#    built from a duplicate of slide 1 (so it inherits the same title /
#    content placeholder geometry) with its picture removed and its text
#    replaced.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Underline part of the run on slide 1 --------------------------------

$s1 = $p.Slides.Item(1)
$body1 = $s1.Shapes.Item(3)
$tr1 = $body1.TextFrame.TextRange

$para2 = $tr1.Paragraphs(2, 1)
$fullText = $para2.Text
$startIdx = $fullText.IndexOf("matched to LSST and CMB-")
$endIdx = $fullText.IndexOf("; Southern site preferable")

$underlineStart = $startIdx + 1
$underlineLength = $endIdx - $startIdx

$toUnderline = $para2.Characters($underlineStart, $underlineLength)
$toUnderline.Font.Underline = $true

# --- 2. Duplicate slide 1 to create the new "Road Map" slide ---------------

$s1.Duplicate() | Out-Null
$s2 = $p.Slides.Item(2)

# Drop the picture that came along with the duplicated slide; the new
# slide only keeps the Title + Content placeholders.
$s2.Shapes.Item(1).Delete()

$title2 = $s2.Shapes.Item(1)
$title2.TextFrame.TextRange.Text = "Road Map for Spectroscopy"

$body2 = $s2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange

$lines = @(
  "The proposed spectroscopic surveys build on each other directly",
  "DESI-2 would be relatively low in cost and could follow DESI immediately in 2024",
  "Spectrograph upgrades to add IR arm would enhance capabilities at higher redshifts",
  "Moving to Blanco is an option, increasing LSST overlap ",
  "SSSI could reuse DESI spectrographs to reduce costs",
  "Earliest possible deployment c. 2026",
  "Most efficient option would be to deploy on 11-12m telescope (e.g. MSE or European wide-field concepts) ",
  "BOA would require both a >10m wide-field telescope and significant hardware R&D",
  "Earliest possible deployment early 2030s",
  "Could utilize telescope originally developed for SSSI",
  "",
  ""
)
$tr2.Text = [string]::Join("`r", $lines)

# Second-level (indented) bullets.
$indentedParas = @(3, 4, 6, 7, 9, 10, 11)
foreach ($idx in $indentedParas) {
  $para = $tr2.Paragraphs($idx, 1)
  $para.IndentLevel = 2
}

# Split "Earliest possible deployment early 2030s" into two runs so the
# second part ("early 2030s") can carry its own (smtClean) run, matching
# how the deck was authored.
$para9 = $tr2.Paragraphs(9, 1)
$splitAt = "Earliest possible deployment ".Length
$firstPart = $para9.Characters(1, $splitAt)
$secondPart = $para9.Characters($splitAt + 1, $para9.Length - $splitAt)
$secondPart.Text = "early 2030s"
